$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "D3"  = -8.7156
    "B4"  = 8.674299999999995
    "D4"  = -7.177699999999996
    "D5"  = -7.503100000000006
    "A8"  = -22.0618
    "A10" = -22.11980000000001
    "A12" = -21.92219999999999
    "B12" = 6.0654
    "C12" = -11.3847
    "C13" = -12.09639999999999
    "B15" = 6.161399999999998
    "B17" = 5.087799999999996
    "A18" = -22.20370000000001
    "D20" = -8.564799999999989
    "C21" = -13.4463
    "D23" = -7.973599999999998
    "C25" = -12.4871
    "B26" = 4.930299999999999
    "D26" = -7.623900000000003
    "B27" = 5.232200000000001
    "B28" = 5.4895
    "C32" = -12.4912
    "D34" = -7.422200000000003
    "C36" = -13.48570000000002
    "A37" = -21.93899999999999
    "B37" = 5.389100000000001
    "C38" = -13.36239999999999
    "D39" = -7.6431
    "D40" = -8.90159999999999
    "C41" = -12.60890000000001
    "D41" = -8.384399999999992
    "B47" = 5.876000000000001
    "D47" = -7.677700000000002
    "C52" = -11.31780000000001
    "A55" = -22.16379999999999
    "C59" = -12.4132
    "D60" = -8.718799999999998
    "B65" = 5.6773
    "C67" = -11.94849999999999
    "A68" = -21.43929999999999
    "D72" = -6.985799999999998
    "B73" = 9.173599999999999
    "A77" = -20.3853
    "A78" = -19.75779999999997
    "A81" = -22.12840000000002
    "A82" = -21.6903
    "D83" = -9.170799999999995
    "B84" = 4.309500000000001
    "C84" = -12.1234
    "B85" = 5.395700000000001
    "C88" = -12.88719999999999
    "C89" = -14.1775
    "B93" = 5.544700000000002
    "B95" = 6.2994
    "C95" = -12.09780000000001
    "B98" = 5.375300000000005
    "B99" = 5.883600000000001
    "B101" = 6.1213
    "C105" = -12.99580000000002
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
